# PA4.docx - "fixed typo in PA description."
#
# This script reproduces the content-level edits from the commit:
#   1. Three places where a sentence had been split across multiple
#      <w:r> runs (an artifact of earlier proofing-tool edits / typo fixes)
#      are re-typed as a single clean run so the paragraph text is one
#      contiguous run again (this also drops the now-stale <w:proofErr/>
#      grammar-check markers that bracketed the old "...>.compressed.txt"
#      split).
#   2. The whole "Possible Strategy for Getting Started" section (a
#      heading + two bullet points) is removed.
#   3. The hidden "_GoBack" bookmark (Word's "last edit position"
#      marker) is moved from the very end of the document to the start
#      of the "Header Comment, and Formatting" heading - i.e. wherever
#      the author's cursor was when the document was last saved.

$d = $word.ActiveDocument

function Set-ParagraphPlainText {
    <#
        Finds the paragraph containing $anchorText and retypes the
        whole paragraph (minus its trailing paragraph mark) as a single
        run containing $newText. Used to collapse a paragraph that is
        split across several runs back into one run.

        (Positional params - named "-param value" binding isn't
        supported by this host's PowerShell subset.)
    #>
    param($anchorText, $newText)

    $r = $d.Content
    $found = $r.Find.Execute($anchorText)
    if (-not $found) {
        Write-Output ("anchor not found: " + $anchorText)
        return
    }

    $r.Expand(4) | Out-Null      # wdParagraph - grow to the full paragraph
    $r.MoveEnd(1, -1) | Out-Null # wdCharacter - back off the paragraph mark
    $r.Delete()
    $r.InsertBefore($newText)
}

# --- 1. Re-merge runs that were split mid-sentence -----------------------

Set-ParagraphPlainText 'Using the algorithm above on the string' 'Using the algorithm above on the string "berries apples berries apples pears apples" (42 chars)'

Set-ParagraphPlainText 'The frequencies of the words are berries' 'The frequencies of the words are berries: 2, apples: 3, pears: 1'

Set-ParagraphPlainText 'To maintain consistency across student submissions' 'To maintain consistency across student submissions, please ensure that the encoded files have the following names: <source>.compressed.txt and <source>.codes.txt.  Thus, if you load in "happy.txt", your program would generate the files "happy.compressed.txt" and "happy.codes.txt."  '

# --- 2. Delete the "Possible Strategy for Getting Started" section -------

$sectionStart = $d.Content
$sectionStart.Find.Execute('Possible Strategy for Getting Started') | Out-Null
$sectionStart.Expand(4) | Out-Null
$startPos = $sectionStart.Start

$sectionEnd = $d.Content
$sectionEnd.Find.Execute('Header Comment, and Formatting') | Out-Null
$sectionEnd.Expand(4) | Out-Null
$endPos = $sectionEnd.Start

$d.Range($startPos, $endPos).Delete()

# --- 3. Move the "_GoBack" bookmark to the "Header Comment" heading ------
# Word keeps one hidden "_GoBack" bookmark marking the place of the last
# edit; re-adding a bookmark with that name moves it (a document can only
# have one bookmark per name), which also removes it from its old spot at
# the end of the document.

$headerComment = $d.Content
$headerComment.Find.Execute('Header Comment, and Formatting') | Out-Null
$headerComment.Collapse(1) | Out-Null
$d.Bookmarks.Add('_GoBack', $headerComment) | Out-Null

Write-Output ("Done. Paragraph count: " + $d.Paragraphs.Count)
